$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 7; this shifts existing rows 7-13 down to 8-14
$ws.Rows.Item(7).Insert()

# Populate the newly inserted row 7. Same constant columns as every other data
# row (A,B,C,E,F,G,H,I,N,O,Q,R), the new date from the diff (44498), and the
# same J/K/L/M/P values that used to be in the old row 7 (now row 8).
$ws.Cells.Item(7, 1).Value = 10
$ws.Cells.Item(7, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(7, 3).Value = 'La Araucanía'
$ws.Cells.Item(7, 4).Value = 44498
$ws.Cells.Item(7, 4).NumberFormat = $ws.Cells.Item(8, 4).NumberFormat
$ws.Cells.Item(7, 5).Value = 9
$ws.Cells.Item(7, 6).Value = 100112017
$ws.Cells.Item(7, 7).Value = 'Ramas de apio'
$ws.Cells.Item(7, 8).Value = 'Sin especificar'
$ws.Cells.Item(7, 9).Value = 'Primera'
$ws.Cells.Item(7, 10).Value = 40
$ws.Cells.Item(7, 11).Value = 4000
$ws.Cells.Item(7, 12).Value = 4000
$ws.Cells.Item(7, 13).Value = 4000
$ws.Cells.Item(7, 14).Value = '$/paquete'
$ws.Cells.Item(7, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(7, 16).Value = 4000
$ws.Cells.Item(7, 17).Value = 1
$ws.Cells.Item(7, 18).Value = 'Hortaliza'
